# Apply cryptos list update (cell text values) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.311.85"
$ws.Range("E2").Value = "  +4.40%  "
$ws.Range("D3").Value = "3.454.64"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "568.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "184.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.42%  "
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("D8").Value = "3.448.26"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +7.83%  "
$ws.Range("E11").Value = "  +2.96%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "55.82"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.99%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000281"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.77%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "9.38"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("D15").Value = "4.007.85"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "3.459.83"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "67.222.46"
$ws.Range("E18").Value = "  +4.33%  "
$ws.Range("E19").Value = "  +1.17%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.04"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("E21").Value = "  +2.76%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "484.70"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.58%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.96"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.88%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "15.06"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +10.85%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "4.18"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "90.08"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("E27").Value = "  +0.37%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.92"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.95"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("E30").Value = "  +3.91%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.97"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.83%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "598.97"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.89%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "11.62"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.96%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "63.13"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("E36").Value = "  +6.53%  "
$ws.Range("E37").Value = "  -0.10%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.63"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "0.0₃0784"
$ws.Range("E39").Value = "  +6.04%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.60"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.388"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.65%  "
$ws.Range("D42").Value = "3.142.44"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("E43").Value = "  +4.40%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.37%  "
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("E46").Value = "  +21.90%  "
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("E48").Value = "  +1.30%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.74"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +7.21%  "
$ws.Range("E50").Value = "  +0.10%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "142.06"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "
